# Updated symbol list on Sun Dec 18 20:40:49 UTC 2022 with GitHub Actions
#
# Applies the per-cell price/volume refresh (and the CEJI <-> BKEXToken
# row-content swap) described by the diff against cryptos.xlsx / Sheet1.
#
# Numeric-looking values are written with a leading apostrophe so Excel
# keeps storing them as text (matching the workbook's existing inlineStr
# "numbers-as-text" convention) instead of silently re-typing the cell as
# a Number; the style is then reset to "Normal" so the quote-prefix flag
# Excel adds for apostrophe-entry doesn't leave a stray cell format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Value
    )
    $rng = $ws.Range($Cell)
    $rng.Value2 = "'" + $Value
    $rng.Style = "Normal"
}

function Set-PlainValue {
    param(
        [string]$Cell,
        [string]$Value
    )
    $ws.Range($Cell).Value2 = $Value
}

# --- Row 2 (BNB) ---
Set-TextValue  "D2"  "249.82"
Set-PlainValue "E2"  "1BNBBNB"

# --- Row 3 (OKB) ---
Set-TextValue "D3" "21.96"

# --- Row 4 (HuobiToken) ---
Set-TextValue "D4" "5.541"

# --- Row 5 (Cronos) ---
Set-TextValue "D5" "0.05646"

# --- Row 6 ---
Set-TextValue "D6" "6.457"

# --- Row 7 ---
Set-TextValue "D7" "0.8010"

# --- Row 8 ---
Set-TextValue "D8" "1.037"

# --- Row 9 ---
Set-TextValue "D9" "0.1439"

# --- Row 10 ---
Set-TextValue "D10" "0.07326"

# --- Row 11 ---
Set-TextValue "D11" "0.03098"

# --- Row 12 ---
Set-TextValue "D12" "0.02915"

# --- Row 13 ---
Set-TextValue "D13" "0.09270"

# --- Row 14 ---
Set-TextValue "D14" "0.001674"

# --- Row 15 ---
Set-TextValue "D15" "3.211"

# --- Row 16 ---
Set-TextValue "D16" "0.04740"

# --- Row 17 (One / ONE) ---
Set-TextValue  "D17" "0.0005811"
Set-PlainValue "E17" "16OneONE"

# --- Row 18 ---
Set-TextValue "D18" "0.006403"

# --- Row 19 ---
Set-TextValue "D19" "0.005070"

# --- Row 20 ---
Set-TextValue "D20" "0.001056"

# --- Row 22 ---
Set-TextValue "D22" "3.976"

# --- Row 23 ---
Set-TextValue "D23" "3.378"

# --- Row 24 ---
Set-TextValue "D24" "2.089"

# --- Row 26 ---
Set-TextValue "D26" "0.1258"

# --- Row 40 (IDEX) ---
Set-TextValue "D40" "0.04152"

# --- Row 41 (KickToken) ---
Set-TextValue "D41" "0.006902"

# --- Row 42: was CEJI, now BKEXToken ---
Set-PlainValue "B42" "BKEXToken"
Set-PlainValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue  "D42" "0.1043"
Set-PlainValue "E42" "41BKEXTokenBKK"

# --- Row 43: was BKEXToken, now CEJI ---
Set-PlainValue "B43" "CEJI"
Set-PlainValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue  "D43" "0.003301"
Set-PlainValue "E43" "42CEJICEJIBestin24h"

# --- Row 44 (LocalTraders) ---
Set-TextValue "D44" "0.009373"

# --- Row 45 (CoinLion) ---
Set-TextValue "D45" "0.00005644"

# --- Row 47 (CoinbaseStockToken) ---
Set-TextValue "D47" "0.6801"

# --- Row 48 (BOLO) ---
Set-TextValue  "D48" "0.01624"
Set-PlainValue "E48" "47BOLOBOLOWorstin24h"
